# "Fuel" sheet update for both CH and SIN
# - update PEN & CO2 for natural gas (NG) in the "FUELS" sheet for SIN
#   (data source: ecoinvent 3.4 market for natural gas, burned in gas motor,
#   for storage, GLO)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")

# PEN (C2): build up the value from its ecoinvent sub-components so the
# formula (and the resulting cached value) are both captured.
$ws.Range("C2").Formula = "=1.1767+0.0019487+0.0000015726"

# CO2 (D2): updated emission factor for natural gas.
$ws.Range("D2").Value = 0.06682

# reference (F2): point at the new ecoinvent source string.
$ws.Range("F2").Value = "ecoinvent 3.4 - market for natural gas, burned in gas motor, for storage_GLO_2017_Allocation, cut-off"

# Leave the cursor where the author last left it when saving.
[void]$ws.Range("C2").Select()
